$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.729584333333333
$ws.Range("H2").Value = 5.188753
$ws.Range("I2").Value = 0.2476387648475193
$ws.Range("J2").Value = 0.2476387648475193
$ws.Range("M2").Value = 4.993165333333334
$ws.Range("N2").Value = 14.979496
$ws.Range("O2").Value = 0.06779298131037136
$ws.Range("P2").Value = 0.06779298131037137
$ws.Range("Q2").Value = 8.636100534276446
$ws.Range("R2").Value = 77.72490480848801
$ws.Range("S2").Value = 0.01678817015703132
$ws.Range("T2").Value = 0.01678817015703133
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.729584333333333
$ws.Range("H3").Value = 5.188753
$ws.Range("I3").Value = 0.2476387648475193
$ws.Range("J3").Value = 0.2476387648475193
$ws.Range("O3").Value = 0.5355771637189464
$ws.Range("P3").Value = 0.5355771637189464
$ws.Range("Q3").Value = 68.22680077401834
$ws.Range("R3").Value = 614.0412069661651
$ws.Range("S3").Value = 0.1326296673038975
$ws.Range("T3").Value = 0.1326296673038975
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.729584333333333
$ws.Range("H4").Value = 5.188753
$ws.Range("I4").Value = 0.2476387648475193
$ws.Range("J4").Value = 0.2476387648475193
$ws.Range("M4").Value = 29.08216166666666
$ws.Range("N4").Value = 87.24648499999999
$ws.Range("O4").Value = 0.3948530262300277
$ws.Range("P4").Value = 0.3948530262300277
$ws.Range("Q4").Value = 50.30005119813389
$ws.Range("R4").Value = 452.700460783205
$ws.Range("S4").Value = 0.09778091571190918
$ws.Range("T4").Value = 0.09778091571190919
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.729584333333333
$ws.Range("H5").Value = 5.188753
$ws.Range("I5").Value = 0.2476387648475193
$ws.Range("J5").Value = 0.2476387648475193
$ws.Range("M5").Value = 0.130869
$ws.Range("N5").Value = 0.392607
$ws.Range("O5").Value = 0.001776828740654623
$ws.Range("P5").Value = 0.001776828740654624
$ws.Range("Q5").Value = 0.226348972119
$ws.Range("R5").Value = 2.037140749071
$ws.Range("S5").Value = 0.0004400116746812841
$ws.Range("T5").Value = 0.0004400116746812842
$ws.Range("I6").Value = 0.2307941364328804
$ws.Range("J6").Value = 0.2307941364328804
$ws.Range("M6").Value = 4.993165333333334
$ws.Range("N6").Value = 14.979496
$ws.Range("O6").Value = 0.06779298131037136
$ws.Range("P6").Value = 0.06779298131037137
$ws.Range("Q6").Value = 8.048664619140446
$ws.Range("R6").Value = 72.43798157226401
$ws.Range("S6").Value = 0.01564622257773756
$ws.Range("T6").Value = 0.01564622257773756
$ws.Range("I7").Value = 0.2307941364328804
$ws.Range("J7").Value = 0.2307941364328804
$ws.Range("O7").Value = 0.5355771637189464
$ws.Range("P7").Value = 0.5355771637189464
$ws.Range("S7").Value = 0.1236080689936856
$ws.Range("T7").Value = 0.1236080689936856
$ws.Range("I8").Value = 0.2307941364328804
$ws.Range("J8").Value = 0.2307941364328804
$ws.Range("M8").Value = 29.08216166666666
$ws.Range("N8").Value = 87.24648499999999
$ws.Range("O8").Value = 0.3948530262300277
$ws.Range("P8").Value = 0.3948530262300277
$ws.Range("Q8").Value = 46.87859304237389
$ws.Range("R8").Value = 421.907337381365
$ws.Range("S8").Value = 0.0911297632066687
$ws.Range("T8").Value = 0.09112976320666871
$ws.Range("I9").Value = 0.2307941364328804
$ws.Range("J9").Value = 0.2307941364328804
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.130869
$ws.Range("N9").Value = 0.392607
$ws.Range("O9").Value = 0.001776828740654623
$ws.Range("P9").Value = 0.001776828740654624
$ws.Range("Q9").Value = 0.210952496007
$ws.Range("R9").Value = 1.898572464063
$ws.Range("S9").Value = 0.0004100816547885062
$ws.Range("T9").Value = 0.0004100816547885063
$ws.Range("G10").Value = 2.743651333333334
$ws.Range("H10").Value = 8.230954000000001
$ws.Range("I10").Value = 0.3928310486309039
$ws.Range("J10").Value = 0.3928310486309038
$ws.Range("M10").Value = 4.993165333333334
$ws.Range("N10").Value = 14.979496
$ws.Range("O10").Value = 0.06779298131037136
$ws.Range("P10").Value = 0.06779298131037137
$ws.Range("Q10").Value = 13.69950472435378
$ws.Range("R10").Value = 123.295542519184
$ws.Range("S10").Value = 0.02663118793796845
$ws.Range("T10").Value = 0.02663118793796845
$ws.Range("G11").Value = 2.743651333333334
$ws.Range("H11").Value = 8.230954000000001
$ws.Range("I11").Value = 0.3928310486309039
$ws.Range("J11").Value = 0.3928310486309038
$ws.Range("O11").Value = 0.5355771637189464
$ws.Range("P11").Value = 0.5355771637189464
$ws.Range("Q11").Value = 108.2286358086634
$ws.Range("R11").Value = 974.0577222779701
$ws.Range("S11").Value = 0.210391338846479
$ws.Range("T11").Value = 0.210391338846479
$ws.Range("G12").Value = 2.743651333333334
$ws.Range("H12").Value = 8.230954000000001
$ws.Range("I12").Value = 0.3928310486309039
$ws.Range("J12").Value = 0.3928310486309038
$ws.Range("M12").Value = 29.08216166666666
$ws.Range("N12").Value = 87.24648499999999
$ws.Range("O12").Value = 0.3948530262300277
$ws.Range("P12").Value = 0.3948530262300277
$ws.Range("Q12").Value = 79.79131163296556
$ws.Range("R12").Value = 718.12180469669
$ws.Range("S12").Value = 0.1551105283490276
$ws.Range("T12").Value = 0.1551105283490276
$ws.Range("G13").Value = 2.743651333333334
$ws.Range("H13").Value = 8.230954000000001
$ws.Range("I13").Value = 0.3928310486309039
$ws.Range("J13").Value = 0.3928310486309038
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.130869
$ws.Range("N13").Value = 0.392607
$ws.Range("O13").Value = 0.001776828740654623
$ws.Range("P13").Value = 0.001776828740654624
$ws.Range("Q13").Value = 0.3590589063420001
$ws.Range("R13").Value = 3.231530157078001
$ws.Range("S13").Value = 0.000697993497428884
$ws.Range("T13").Value = 0.000697993497428884
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.8991316666666668
$ws.Range("H14").Value = 2.697395
$ws.Range("I14").Value = 0.1287360500886965
$ws.Range("J14").Value = 0.1287360500886965
$ws.Range("M14").Value = 4.993165333333334
$ws.Range("N14").Value = 14.979496
$ws.Range("O14").Value = 0.06779298131037136
$ws.Range("P14").Value = 0.06779298131037137
$ws.Range("Q14").Value = 4.489513068102223
$ws.Range("R14").Value = 40.40561761292
$ws.Range("S14").Value = 0.008727400637634034
$ws.Range("T14").Value = 0.008727400637634036
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.8991316666666668
$ws.Range("H15").Value = 2.697395
$ws.Range("I15").Value = 0.1287360500886965
$ws.Range("J15").Value = 0.1287360500886965
$ws.Range("O15").Value = 0.5355771637189464
$ws.Range("P15").Value = 0.5355771637189464
$ws.Range("Q15").Value = 35.46798841144167
$ws.Range("R15").Value = 319.211895702975
$ws.Range("S15").Value = 0.0689480885748843
$ws.Range("T15").Value = 0.0689480885748843
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.8991316666666668
$ws.Range("H16").Value = 2.697395
$ws.Range("I16").Value = 0.1287360500886965
$ws.Range("J16").Value = 0.1287360500886965
$ws.Range("M16").Value = 29.08216166666666
$ws.Range("N16").Value = 87.24648499999999
$ws.Range("O16").Value = 0.3948530262300277
$ws.Range("P16").Value = 0.3948530262300277
$ws.Range("Q16").Value = 26.14869248961945
$ws.Range("R16").Value = 235.338232406575
$ws.Range("S16").Value = 0.05083181896242224
$ws.Range("T16").Value = 0.05083181896242225
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.8991316666666668
$ws.Range("H17").Value = 2.697395
$ws.Range("I17").Value = 0.1287360500886965
$ws.Range("J17").Value = 0.1287360500886965
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.130869
$ws.Range("N17").Value = 0.392607
$ws.Range("O17").Value = 0.001776828740654623
$ws.Range("P17").Value = 0.001776828740654624
$ws.Range("Q17").Value = 0.117668462085
$ws.Range("R17").Value = 1.059016158765
$ws.Range("S17").Value = 0.0002287419137559491
$ws.Range("T17").Value = 0.0002287419137559492
